# Publica dados no site (31/10/2025 10:41)
# Adds the new "UFPA" competition-ratio sheet, positioned right before the
# existing "Enare" sheet (same place it was inserted in the source commit).

$wb = $excel.ActiveWorkbook

# Sheet that already carries the "PROGRAMA / TOTAL DE VAGAS / INSCRITOS /
# CONCORRENCIA" header look (bold, centered, boxed) - used purely as a
# formatting donor for the new sheet's header row.
$headerDonor = $wb.Worksheets.Item("Instituto do Câncer do Ceará – ")

# Anchor sheet: the new sheet must land immediately before this one.
$enare = $wb.Worksheets.Item("Enare")

# Insert the new sheet right before "Enare" and name it.
$ws = $wb.Worksheets.Add($enare)
$ws.Name = "UFPA"

# --- Header row (row 1) -----------------------------------------------
$headerDonor.Range("A1:D1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(1, 1).Value = "PROGRAMA"
$ws.Cells.Item(1, 2).Value = "VAGAS"
$ws.Cells.Item(1, 3).Value = "INSCRITOS"
$ws.Cells.Item(1, 4).Value = "CONCORRÊNCIA"

# --- Data rows (rows 2-16) ---------------------------------------------
# PROGRAMA is stored as plain text (never numeric-looking, no special
# handling needed). VAGAS / INSCRITOS / CONCORRENCIA look numeric but must
# be stored as text (matching every other sheet in this workbook), so they
# are entered with a leading apostrophe to force text storage.
$data = @(
    @("Clínica Médica (Belém)", "10", "200", "20"),
    @("Clínica Médica (Altamira)", "2", "6", "3"),
    @("Cirurgia Geral (Belém)", "6", "169", "28,17"),
    @("Dermatologia (Belém)", "3", "70", "23,33"),
    @("Ginecologia e Obstetrícia (Belém)", "6", "95", "15,83"),
    @("Medicina de Família e Comunidade (Belém)", "10", "133", "13,3"),
    @("Medicina de Família e Comunidade (Altamira)", "2", "4", "2"),
    @("Oftalmologia (Belém)", "4", "88", "22"),
    @("Otorrinolaringologia (Belém)", "4", "90", "22,5"),
    @("Pediatria (Belém)", "8", "92", "11,5"),
    @("Endocrinologia (Belém)", "2", "20", "10"),
    @("Geriatria (Belém)", "2", "3", "1,5"),
    @("Oncologia Clínica (Belém)", "1", "7", "7"),
    @("Pneumologia (Belém)", "2", "9", "4,5"),
    @("Medicina Paliativa (Belém)", "2", "4", "2")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = "'" + $r[1]
    $ws.Cells.Item($row, 3).Value = "'" + $r[2]
    $ws.Cells.Item($row, 4).Value = "'" + $r[3]
    $row = $row + 1
}
